# Portfolio workbook refresh: new price/RSI snapshot (watchlist re-sorted by
# RSI ascending), corresponding stocks-sheet recompute, updated portfolio
# totals, and a new summary row appended for the latest snapshot.
# (commit: "comment GS for temp fix")

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("watchlist")
$ws2 = $wb.Worksheets.Item("stocks")
$ws3 = $wb.Worksheets.Item("portfolio")
$ws5 = $wb.Worksheets.Item("summary")

# ---------------------------------------------------------------------------
# watchlist: ticker / price / pct_change / rsi, rows 2-31, re-sorted by rsi asc
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = "VZ"
$ws1.Range("B2").Value = 54.1
$ws1.Range("C2").Value = 0.2409
$ws1.Range("D2").Value = 22.64150943396223

$ws1.Range("A3").Value = "JNJ"
$ws1.Range("B3").Value = 144.37
$ws1.Range("C3").Value = -0.9127999999999999
$ws1.Range("D3").Value = 34.71837488457984

$ws1.Range("A4").Value = "PG"
$ws1.Range("B4").Value = 112.6
$ws1.Range("C4").Value = 0.878
$ws1.Range("D4").Value = 35.79304495335023

$ws1.Range("A5").Value = "IBM"
$ws1.Range("B5").Value = 118.39
$ws1.Range("C5").Value = -0.6128
$ws1.Range("D5").Value = 37.52711496746203

$ws1.Range("A6").Value = "WBA"
$ws1.Range("B6").Value = 39.6
$ws1.Range("C6").Value = -0.3774
$ws1.Range("D6").Value = 40.31476997578689

$ws1.Range("A7").Value = "CVX"
$ws1.Range("B7").Value = 90.28
$ws1.Range("C7").Value = -1.9122
$ws1.Range("D7").Value = 41.56706507304116

$ws1.Range("A8").Value = "KO"
$ws1.Range("B8").Value = 45.03
$ws1.Range("C8").Value = -0.3099
$ws1.Range("D8").Value = 42.50000000000002

$ws1.Range("A9").Value = "GS"
$ws1.Range("B9").Value = 179.93
$ws1.Range("C9").Value = -0.0944
$ws1.Range("D9").Value = 42.81578947368421

$ws1.Range("A10").Value = "JPM"
$ws1.Range("B10").Value = 89.47
$ws1.Range("C10").Value = -0.7763
$ws1.Range("D10").Value = 43.44295574502638

$ws1.Range("A11").Value = "XOM"
$ws1.Range("B11").Value = 44.6
$ws1.Range("C11").Value = 0.0898
$ws1.Range("D11").Value = 46.83544303797471

$ws1.Range("A12").Value = "MSFT"
$ws1.Range("B12").Value = 183.51
$ws1.Range("C12").Value = 0.0436
$ws1.Range("D12").Value = 48.37438423645321

$ws1.Range("A13").Value = "MMM"
$ws1.Range("B13").Value = 146.44
$ws1.Range("C13").Value = 0.4321
$ws1.Range("D13").Value = 49.00813008130079

$ws1.Range("A14").Value = "AXP"
$ws1.Range("B14").Value = 89.33
$ws1.Range("C14").Value = -0.5566
$ws1.Range("D14").Value = 49.72624798711755

$ws1.Range("A15").Value = "MRK"
$ws1.Range("B15").Value = 76.37
$ws1.Range("C15").Value = -0.2351
$ws1.Range("D15").Value = 49.82896237172177

$ws1.Range("A16").Value = "CAT"
$ws1.Range("B16").Value = 112.47
$ws1.Range("C16").Value = -1.394
$ws1.Range("D16").Value = 50.64423765211167

$ws1.Range("A17").Value = "PFE"
$ws1.Range("B17").Value = 37.5
$ws1.Range("C17").Value = 0.6441
$ws1.Range("D17").Value = 53.24074074074074

$ws1.Range("A18").Value = "UNH"
$ws1.Range("B18").Value = 289.94
$ws1.Range("C18").Value = 1.0561
$ws1.Range("D18").Value = 53.94948952176249

$ws1.Range("A19").Value = "BA"
$ws1.Range("B19").Value = 137.53
$ws1.Range("C19").Value = -1.0576
$ws1.Range("D19").Value = 54.45048966267682

$ws1.Range("A20").Value = "TRV"
$ws1.Range("B20").Value = 100.1
$ws1.Range("C20").Value = 0.9378
$ws1.Range("D20").Value = 56.33608815426995

$ws1.Range("A21").Value = "MCD"
$ws1.Range("B21").Value = 184.41
$ws1.Range("C21").Value = -0.362
$ws1.Range("D21").Value = 56.34477254588986

$ws1.Range("A22").Value = "RTX"
$ws1.Range("B22").Value = 60
$ws1.Range("C22").Value = -0.1664
$ws1.Range("D22").Value = 56.77655677655678

$ws1.Range("A23").Value = "WMT"
$ws1.Range("B23").Value = 124.33
$ws1.Range("C23").Value = -0.528
$ws1.Range("D23").Value = 57.11361310133061

$ws1.Range("A24").Value = "V"
$ws1.Range("B24").Value = 190.86
$ws1.Range("C24").Value = 0.1259
$ws1.Range("D24").Value = 59.50594121325831

$ws1.Range("A25").Value = "DOW"
$ws1.Range("B25").Value = 36.12
$ws1.Range("C25").Value = 0.8375
$ws1.Range("D25").Value = 60.32295271049597

$ws1.Range("A26").Value = "HD"
$ws1.Range("B26").Value = 241.88
$ws1.Range("C26").Value = 0.4151
$ws1.Range("D26").Value = 60.68253513048464

$ws1.Range("A27").Value = "INTC"
$ws1.Range("B27").Value = 62.26
$ws1.Range("C27").Value = 0.4518
$ws1.Range("D27").Value = 61.55218554861729

$ws1.Range("A28").Value = "NKE"
$ws1.Range("B28").Value = 93.75
$ws1.Range("C28").Value = -0.5411
$ws1.Range("D28").Value = 61.64189667374379

$ws1.Range("A29").Value = "AAPL"
$ws1.Range("B29").Value = 318.89
$ws1.Range("C29").Value = 0.6438
$ws1.Range("D29").Value = 62.32414181204273

$ws1.Range("A30").Value = "CSCO"
$ws1.Range("B30").Value = 44.9
$ws1.Range("C30").Value = 0.5824
$ws1.Range("D30").Value = 63.35664335664335

$ws1.Range("A31").Value = "DIS"
$ws1.Range("B31").Value = 118.02
$ws1.Range("C31").Value = 0.1612
$ws1.Range("D31").Value = 64.53412073490814

# ---------------------------------------------------------------------------
# stocks: current_price / value / performance / current_rsi refresh, rows 2-20
# (ticker / purch_price / shares / last_activity are unchanged)
# ---------------------------------------------------------------------------
$ws2.Range("C2").Value = 100.1
$ws2.Range("E2").Value = 1001
$ws2.Range("F2").Value = 4.2166
$ws2.Range("G2").Value = 56.33608815426995

$ws2.Range("C3").Value = 39.6
$ws2.Range("E3").Value = 950.4000000000001
$ws2.Range("F3").Value = -4.9448
$ws2.Range("G3").Value = 40.31476997578689

$ws2.Range("C4").Value = 112.47
$ws2.Range("E4").Value = 787.29
$ws2.Range("F4").Value = 2.0599
$ws2.Range("G4").Value = 50.64423765211167

$ws2.Range("C5").Value = 60
$ws2.Range("E5").Value = 720
$ws2.Range("F5").Value = 1.6088
$ws2.Range("G5").Value = 56.77655677655678

$ws2.Range("C6").Value = 54.1
$ws2.Range("E6").Value = 595.1
$ws2.Range("F6").Value = -4.3324
$ws2.Range("G6").Value = 22.64150943396223

$ws2.Range("C7").Value = 124.33
$ws2.Range("E7").Value = 372.99
$ws2.Range("F7").Value = -0.008
$ws2.Range("G7").Value = 57.11361310133061

$ws2.Range("C8").Value = 112.6
$ws2.Range("E8").Value = 450.4
$ws2.Range("F8").Value = -1.7452
$ws2.Range("G8").Value = 35.79304495335023

$ws2.Range("C9").Value = 76.37
$ws2.Range("E9").Value = 458.22
$ws2.Range("F9").Value = -1.9766
$ws2.Range("G9").Value = 49.82896237172177

$ws2.Range("C10").Value = 144.37
$ws2.Range("E10").Value = 433.11
$ws2.Range("F10").Value = -3.1789
$ws2.Range("G10").Value = 34.71837488457984

$ws2.Range("C11").Value = 146.44
$ws2.Range("E11").Value = 292.88
$ws2.Range("F11").Value = -0.041
$ws2.Range("G11").Value = 49.00813008130079

$ws2.Range("C12").Value = 89.47
$ws2.Range("E12").Value = 357.88
$ws2.Range("F12").Value = 1.4514
$ws2.Range("G12").Value = 43.44295574502638

$ws2.Range("C13").Value = 184.41
$ws2.Range("E13").Value = 184.41
$ws2.Range("F13").Value = 4.4579
$ws2.Range("G13").Value = 56.34477254588986

$ws2.Range("C14").Value = 36.12
$ws2.Range("E14").Value = 325.08
$ws2.Range("F14").Value = 9.1568
$ws2.Range("G14").Value = 60.32295271049597

$ws2.Range("C15").Value = 118.39
$ws2.Range("E15").Value = 236.78
$ws2.Range("F15").Value = -0.4289
$ws2.Range("G15").Value = 37.52711496746203

$ws2.Range("C16").Value = 45.03
$ws2.Range("E16").Value = 270.18
$ws2.Range("F16").Value = 1.6249
$ws2.Range("G16").Value = 42.50000000000002

$ws2.Range("C17").Value = 89.33
$ws2.Range("E17").Value = 267.99
$ws2.Range("F17").Value = 10.7488
$ws2.Range("G17").Value = 49.72624798711755

$ws2.Range("C18").Value = 179.93
$ws2.Range("E18").Value = 179.93
$ws2.Range("F18").Value = 4.5558
$ws2.Range("G18").Value = 42.81578947368421

$ws2.Range("C19").Value = 137.53
$ws2.Range("E19").Value = 137.53
$ws2.Range("F19").Value = 13.8399
$ws2.Range("G19").Value = 54.45048966267682

$ws2.Range("C20").Value = 62.26
$ws2.Range("E20").Value = 186.78
$ws2.Range("F20").Value = 8.4102
$ws2.Range("G20").Value = 61.55218554861729

# ---------------------------------------------------------------------------
# portfolio: STOCKS / TOTAL values refresh
# ---------------------------------------------------------------------------
$ws3.Range("B3").Value = 8207.950000000001
$ws3.Range("B4").Value = 10066.49

# ---------------------------------------------------------------------------
# summary: append new snapshot row 15 (carry the row-14 date-cell formatting
# over to A15, same as every other row in the column, then fill the values)
# ---------------------------------------------------------------------------
$ws5.Range("A14").Copy()
$ws5.Range("A15").PasteSpecial(-4122)

$ws5.Range("A15").Value = "25/05/2020 07:36:29"
$ws5.Range("B15").Value = 1858.54
$ws5.Range("C15").Value = 8207.950000000001
$ws5.Range("D15").Value = 10066.49

Write-Host "applied portfolio refresh"
